$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.321822762489319
$ws.Range("B1").Value = 1.379718065261841
$ws.Range("C1").Value = 1.551220536231995
$ws.Range("D1").Value = 2.397186040878296
$ws.Range("E1").Value = -1
